$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
